# Apply the "term 2.0" update:
#  - Metadata sheet: bump Version, Date, Contact
#  - Add a new worksheet "Include from FSIII 2" that mirrors
#    "Include from FSIII" (same Property/Operation/Value table)

$wb = $excel.ActiveWorkbook

# --- 1. Update the Metadata sheet -----------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-03T10:45:43+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# --- 2. Add the new "Include from FSIII 2" sheet --------------------------
# Clone the existing "Include from FSIII" sheet (same table layout, styles
# and column widths) straight after itself, then rename it.
$src = $wb.Worksheets.Item("Include from FSIII")
$src.Copy($null, $src)
$ws = $wb.Worksheets.Item($src.Index + 1)
$ws.Name = "Include from FSIII 2"
